# Fruta / hortaliza, semanal
# The data rows (2..54) get their (Fecha, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Precio $/Kg) tuples reassigned according to a
# permutation of rows - i.e. row r receives the values that used to live in
# row $map[r]. Row 23 is a fixed point (keeps its own values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    2=40; 3=7; 4=26; 5=30; 6=18; 7=25; 8=2; 9=5; 10=14; 11=32; 12=50; 13=16;
    14=37; 15=36; 16=49; 17=11; 18=43; 19=13; 20=53; 21=27; 22=6; 23=23;
    24=22; 25=17; 26=51; 27=31; 28=9; 29=28; 30=12; 31=45; 32=47; 33=54;
    34=44; 35=33; 36=42; 37=35; 38=24; 39=52; 40=21; 41=39; 42=29; 43=46;
    44=15; 45=3; 46=38; 47=41; 48=10; 49=34; 50=20; 51=8; 52=4; 53=19; 54=48
}

$cols = @("D", "M", "N", "O", "P", "S")

# Snapshot all the original values first, since we will overwrite them in place.
$original = @{}
foreach ($row in 2..54) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $original[$row] = $rowVals
}

foreach ($row in 2..54) {
    $srcRow = $map[$row]
    $srcVals = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $srcVals[$col]
    }
}
